$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B <-> C (values) and D <-> E (values), row by row,
# mirroring the user having dragged column C before B, then E before D.

1..4 | ForEach-Object {
    $r = $_
    $colB = $ws.Cells.Item($r, 2)
    $colC = $ws.Cells.Item($r, 3)
    $colD = $ws.Cells.Item($r, 4)
    $colE = $ws.Cells.Item($r, 5)

    $bVal = $colB.Value2
    $cVal = $colC.Value2
    $dVal = $colD.Value2
    $eVal = $colE.Value2

    $colB.Value = $cVal
    $colC.Value = $bVal
    $colD.Value = $eVal
    $colE.Value = $dVal
}

# Swap the column widths to match (B<->C, D<->E)
$wB = $ws.Columns.Item(2).ColumnWidth
$wC = $ws.Columns.Item(3).ColumnWidth
$wD = $ws.Columns.Item(4).ColumnWidth
$wE = $ws.Columns.Item(5).ColumnWidth

$ws.Columns.Item(2).ColumnWidth = $wC
$ws.Columns.Item(3).ColumnWidth = $wB
$ws.Columns.Item(4).ColumnWidth = $wE
$ws.Columns.Item(5).ColumnWidth = $wD

# Final selection lands on entire column D, matching the last move (E -> D)
$ws.Columns.Item(4).Select()
